$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 13.40304848865214
$ws.Cells.Item(2, 3).Value = 10.15474694556152
$ws.Cells.Item(2, 4).Value = 5.63351269681556
$ws.Cells.Item(2, 6).Value = 27.43045188391782
$ws.Cells.Item(2, 7).Value = 3.653897818896758
$ws.Cells.Item(2, 11).Value = 9.359707765906203
$ws.Cells.Item(2, 12).Value = 10.75192066467547
$ws.Cells.Item(2, 14).Value = 19.57132153627876
$ws.Cells.Item(2, 15).Value = 24.58727811604889

$ws.Cells.Item(3, 2).Value = 13.15482871798562
$ws.Cells.Item(3, 3).Value = 10.17657629031988
$ws.Cells.Item(3, 4).Value = 5.58942561933071
$ws.Cells.Item(3, 6).Value = 27.44966090157199
$ws.Cells.Item(3, 7).Value = 3.655713877295331
$ws.Cells.Item(3, 11).Value = 9.187848480459458
$ws.Cells.Item(3, 12).Value = 10.72587913264129
$ws.Cells.Item(3, 14).Value = 19.63263712833589
$ws.Cells.Item(3, 15).Value = 24.64654678926594

$ws.Cells.Item(4, 2).Value = 13.00284722521983
$ws.Cells.Item(4, 3).Value = 10.19070000833979
$ws.Cells.Item(4, 4).Value = 5.561701747500542
$ws.Cells.Item(4, 6).Value = 27.46808610621336
$ws.Cells.Item(4, 7).Value = 3.656888737965986
$ws.Cells.Item(4, 11).Value = 9.082335904715855
$ws.Cells.Item(4, 12).Value = 10.71203730312298
$ws.Cells.Item(4, 14).Value = 19.67203517872317
$ws.Cells.Item(4, 15).Value = 24.6878315824622

$ws.Cells.Item(5, 2).Value = 12.94110751347257
$ws.Cells.Item(5, 3).Value = 10.19663723886631
$ws.Cells.Item(5, 4).Value = 5.550243408841828
$ws.Cells.Item(5, 6).Value = 27.47726069683019
$ws.Cells.Item(5, 7).Value = 3.657382585419827
$ws.Cells.Item(5, 11).Value = 9.039398974285429
$ws.Cells.Item(5, 12).Value = 10.70694091929049
$ws.Cells.Item(5, 14).Value = 19.68853155759626
$ws.Cells.Item(5, 15).Value = 24.70588401717098

$ws.Cells.Item(6, 2).Value = 12.93086995345059
$ws.Cells.Item(6, 3).Value = 10.19763410187626
$ws.Cells.Item(6, 4).Value = 5.548331172487697
$ws.Cells.Item(6, 6).Value = 27.47888471257365
$ws.Cells.Item(6, 7).Value = 3.657465500712636
$ws.Cells.Item(6, 11).Value = 9.032274679940297
$ws.Cells.Item(6, 12).Value = 10.7061276542787
$ws.Cells.Item(6, 14).Value = 19.69129746348526
$ws.Cells.Item(6, 15).Value = 24.70895575006997

$ws.Cells.Item(7, 2).Value = 13.00201368103877
$ws.Cells.Item(7, 3).Value = 10.19077934347943
$ws.Cells.Item(7, 4).Value = 5.561547861605812
$ws.Cells.Item(7, 6).Value = 27.46820309388101
$ws.Cells.Item(7, 7).Value = 3.656895337036248
$ws.Cells.Item(7, 11).Value = 9.081756520997653
$ws.Cells.Item(7, 12).Value = 10.71196636258738
$ws.Cells.Item(7, 14).Value = 19.67225586592306
$ws.Cells.Item(7, 15).Value = 24.68807007233613

$ws.Cells.Item(8, 2).Value = 13.31742545540147
$ws.Cells.Item(8, 3).Value = 10.16212451178198
$ws.Cells.Item(8, 4).Value = 5.618449564192366
$ws.Cells.Item(8, 6).Value = 27.43569844565939
$ws.Cells.Item(8, 7).Value = 3.654511613812591
$ws.Cells.Item(8, 11).Value = 9.300483066912001
$ws.Cells.Item(8, 12).Value = 10.74249844005322
$ws.Cells.Item(8, 14).Value = 19.59210071694698
$ws.Cells.Item(8, 15).Value = 24.60669714367386

$ws.Cells.Item(9, 2).Value = 13.93559563981092
$ws.Cells.Item(9, 3).Value = 10.11162368905159
$ws.Cells.Item(9, 4).Value = 5.724676601732117
$ws.Cells.Item(9, 6).Value = 27.42459073246336
$ws.Cells.Item(9, 7).Value = 3.650309446515106
$ws.Cells.Item(9, 11).Value = 9.72699900225985
$ws.Cells.Item(9, 12).Value = 10.81920745208661
$ws.Cells.Item(9, 14).Value = 19.44874223798072
$ws.Cells.Item(9, 15).Value = 24.48603680635857

$ws.Cells.Item(10, 2).Value = 14.38470523143157
$ws.Cells.Item(10, 3).Value = 10.07795541274686
$ws.Cells.Item(10, 4).Value = 5.799239756431397
$ws.Cells.Item(10, 6).Value = 27.44849895885827
$ws.Cells.Item(10, 7).Value = 3.647507067345872
$ws.Cells.Item(10, 11).Value = 10.03568020021216
$ws.Cells.Item(10, 12).Value = 10.88552252133271
$ws.Cells.Item(10, 14).Value = 19.35175938088674
$ws.Cells.Item(10, 15).Value = 24.42121480222752

$ws.Cells.Item(11, 2).Value = 14.58695048162941
$ws.Cells.Item(11, 3).Value = 10.06337726589893
$ws.Cells.Item(11, 4).Value = 5.832359869051167
$ws.Cells.Item(11, 6).Value = 27.4663172365292
$ws.Cells.Item(11, 7).Value = 3.646293434383687
$ws.Cells.Item(11, 11).Value = 10.17445120569772
$ws.Cells.Item(11, 12).Value = 10.91777784508146
$ws.Cells.Item(11, 14).Value = 19.30943256395931
$ws.Cells.Item(11, 15).Value = 24.39691824467169

$ws.Cells.Item(12, 2).Value = 14.66316392241552
$ws.Cells.Item(12, 3).Value = 10.05796242817156
$ws.Cells.Item(12, 4).Value = 5.844782976121117
$ws.Cells.Item(12, 6).Value = 27.47405943378238
$ws.Cells.Item(12, 7).Value = 3.6458426133389
$ws.Cells.Item(12, 11).Value = 10.22671290635664
$ws.Cells.Item(12, 12).Value = 10.93028530347981
$ws.Cells.Item(12, 14).Value = 19.29366071749289
$ws.Cells.Item(12, 15).Value = 24.38846545230449

$ws.Cells.Item(13, 2).Value = 14.64676773620258
$ws.Cells.Item(13, 3).Value = 10.05912392201104
$ws.Cells.Item(13, 4).Value = 5.842112780223522
$ws.Cells.Item(13, 6).Value = 27.47234783441739
$ws.Cells.Item(13, 7).Value = 3.645939316993043
$ws.Cells.Item(13, 11).Value = 10.21547099455314
$ws.Cells.Item(13, 12).Value = 10.92757868380398
$ws.Cells.Item(13, 14).Value = 19.29704608182309
$ws.Cells.Item(13, 15).Value = 24.39025263631923

$ws.Cells.Item(14, 2).Value = 14.59322849627674
$ws.Cells.Item(14, 3).Value = 10.06292967056412
$ws.Cells.Item(14, 4).Value = 5.833384327503129
$ws.Cells.Item(14, 6).Value = 27.46693427772683
$ws.Cells.Item(14, 7).Value = 3.646256169836421
$ws.Cells.Item(14, 11).Value = 10.17875685701728
$ws.Cells.Item(14, 12).Value = 10.91880101058922
$ws.Cells.Item(14, 14).Value = 19.30812987266204
$ws.Cells.Item(14, 15).Value = 24.39620783504411

$ws.Cells.Item(15, 2).Value = 14.56038345658408
$ws.Cells.Item(15, 3).Value = 10.06527453783938
$ws.Cells.Item(15, 4).Value = 5.82802231888276
$ws.Cells.Item(15, 6).Value = 27.46374775113585
$ws.Cells.Item(15, 7).Value = 3.646451390304954
$ws.Cells.Item(15, 11).Value = 10.15622944188333
$ws.Cells.Item(15, 12).Value = 10.91346237666644
$ws.Cells.Item(15, 14).Value = 19.31495237102502
$ws.Cells.Item(15, 15).Value = 24.39995298438459

$ws.Cells.Item(16, 2).Value = 14.37143999259275
$ws.Cells.Item(16, 3).Value = 10.07892293147174
$ws.Cells.Item(16, 4).Value = 5.797058868636726
$ws.Cells.Item(16, 6).Value = 27.44747391193988
$ws.Cells.Item(16, 7).Value = 3.647587608623745
$ws.Cells.Item(16, 11).Value = 10.02657360243156
$ws.Cells.Item(16, 12).Value = 10.88345597373222
$ws.Cells.Item(16, 14).Value = 19.35456148079368
$ws.Cells.Item(16, 15).Value = 24.42290721650119

$ws.Cells.Item(17, 2).Value = 14.25494665021858
$ws.Cells.Item(17, 3).Value = 10.08748437427322
$ws.Cells.Item(17, 4).Value = 5.777856369488048
$ws.Cells.Item(17, 6).Value = 27.43926633753497
$ws.Cells.Item(17, 7).Value = 3.648300281811875
$ws.Cells.Item(17, 11).Value = 9.946574545917107
$ws.Cells.Item(17, 12).Value = 10.86557787118115
$ws.Cells.Item(17, 14).Value = 19.37931831436364
$ws.Cells.Item(17, 15).Value = 24.43831949509327

$ws.Cells.Item(18, 2).Value = 14.18775348871223
$ws.Cells.Item(18, 3).Value = 10.09247815666603
$ws.Cells.Item(18, 4).Value = 5.766736616239563
$ws.Cells.Item(18, 6).Value = 27.43519926796891
$ws.Cells.Item(18, 7).Value = 3.648715953978212
$ws.Cells.Item(18, 11).Value = 9.90040879806104
$ws.Cells.Item(18, 12).Value = 10.85549187735156
$ws.Cells.Item(18, 14).Value = 19.39372645284874
$ws.Cells.Item(18, 15).Value = 24.44767279669155

$ws.Cells.Item(19, 2).Value = 14.16497293260931
$ws.Cells.Item(19, 3).Value = 10.09418091171098
$ws.Cells.Item(19, 4).Value = 5.762958899287489
$ws.Cells.Item(19, 6).Value = 27.43393459477689
$ws.Cells.Item(19, 7).Value = 3.648857684263124
$ws.Cells.Item(19, 11).Value = 9.884753280613266
$ws.Cells.Item(19, 12).Value = 10.85211098414707
$ws.Cells.Item(19, 14).Value = 19.39863380470976
$ws.Cells.Item(19, 15).Value = 24.45092353290284

$ws.Cells.Item(20, 2).Value = 14.26736772243932
$ws.Cells.Item(20, 3).Value = 10.0865658080932
$ws.Cells.Item(20, 4).Value = 5.779908297039267
$ws.Cells.Item(20, 6).Value = 27.44007241365029
$ws.Cells.Item(20, 7).Value = 3.648223820561594
$ws.Cells.Item(20, 11).Value = 9.955106731050222
$ws.Cells.Item(20, 12).Value = 10.86746068606417
$ws.Cells.Item(20, 14).Value = 19.37666545970133
$ws.Cells.Item(20, 15).Value = 24.43662825688364

$ws.Cells.Item(21, 2).Value = 14.60896496769313
$ws.Cells.Item(21, 3).Value = 10.06180896818193
$ws.Cells.Item(21, 4).Value = 5.835951339129705
$ws.Cells.Item(21, 6).Value = 27.46849740431871
$ws.Cells.Item(21, 7).Value = 3.64616286519432
$ws.Cells.Item(21, 11).Value = 10.18954889200009
$ws.Cells.Item(21, 12).Value = 10.92137133189515
$ws.Cells.Item(21, 14).Value = 19.30486734523511
$ws.Cells.Item(21, 15).Value = 24.39443834498963

$ws.Cells.Item(22, 2).Value = 14.83001496757477
$ws.Cells.Item(22, 3).Value = 10.04624417686929
$ws.Cells.Item(22, 4).Value = 5.87188474742976
$ws.Cells.Item(22, 6).Value = 27.49287097575355
$ws.Cells.Item(22, 7).Value = 3.644866925449532
$ws.Cells.Item(22, 11).Value = 10.34107093506111
$ws.Cells.Item(22, 12).Value = 10.95830963094233
$ws.Cells.Item(22, 14).Value = 19.25943709088216
$ws.Cells.Item(22, 15).Value = 24.3712239065874

$ws.Cells.Item(23, 2).Value = 14.7122623122536
$ws.Cells.Item(23, 3).Value = 10.05449526652247
$ws.Cells.Item(23, 4).Value = 5.852771161934752
$ws.Cells.Item(23, 6).Value = 27.47933338197285
$ws.Cells.Item(23, 7).Value = 3.645553939129994
$ws.Cells.Item(23, 11).Value = 10.26037229542625
$ws.Cells.Item(23, 12).Value = 10.9384414916237
$ws.Cells.Item(23, 14).Value = 19.28354775790563
$ws.Cells.Item(23, 15).Value = 24.38321466496112

$ws.Cells.Item(24, 2).Value = 14.26175283376709
$ws.Cells.Item(24, 3).Value = 10.0869808684991
$ws.Cells.Item(24, 4).Value = 5.778980869418217
$ws.Cells.Item(24, 6).Value = 27.43970595676584
$ws.Cells.Item(24, 7).Value = 3.648258370165567
$ws.Cells.Item(24, 11).Value = 9.951249865818095
$ws.Cells.Item(24, 12).Value = 10.86660886572118
$ws.Cells.Item(24, 14).Value = 19.37786426974693
$ws.Cells.Item(24, 15).Value = 24.43739133141372

$ws.Cells.Item(25, 2).Value = 13.76891292759486
$ws.Cells.Item(25, 3).Value = 10.12467983683945
$ws.Cells.Item(25, 4).Value = 5.696535413416036
$ws.Cells.Item(25, 6).Value = 27.42195923822763
$ws.Cells.Item(25, 7).Value = 3.651395989547849
$ws.Cells.Item(25, 11).Value = 9.612212010965418
$ws.Cells.Item(25, 12).Value = 10.79668296756993
$ws.Cells.Item(25, 14).Value = 19.48605329493356
$ws.Cells.Item(25, 15).Value = 24.51450119554295
